$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.950.84'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.884.74'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('D4').Formula = "'1.000"
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Formula = "'305.55"
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Formula = "'0.9998"
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Formula = "'0.5160"
$ws.Range('E7').Value = '  +2.46%  '
$ws.Range('D8').Formula = "'0.3742"
$ws.Range('E8').Value = '  +2.66%  '
$ws.Range('D9').Formula = "'0.07183"
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').Formula = "'21.06"
$ws.Range('E10').Value = '  +1.76%  '
$ws.Range('D11').Formula = "'0.8994"
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').Formula = "'0.07651"
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('D13').Value = '1.858.55'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').Formula = "'93.70"
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Formula = "'5.230"
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Formula = "'0.000008465"
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').Formula = "'14.40"
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').Formula = "'0.9993"
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '26.990.48'
$ws.Range('D21').Formula = "'5.034"
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '2.099.84'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').Formula = "'10.55"
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Formula = "'6.377"
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Formula = "'2.297"
$ws.Range('E25').Value = '  +10.23%  '
$ws.Range('D26').Formula = "'146.31"
$ws.Range('E26').Value = '  -1.08%  '
$ws.Range('D27').Formula = "'18.01"
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').Formula = "'1.726"
$ws.Range('E28').Value = '  -2.89%  '
$ws.Range('D29').Formula = "'113.90"
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('D30').Formula = "'4.914"
$ws.Range('E30').Value = '  +5.27%  '
$ws.Range('D31').Formula = "'4.774"
$ws.Range('E31').Value = '  +1.77%  '
$ws.Range('D32').Formula = "'0.09178"
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('D33').Formula = "'0.05025"
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Formula = "'1.231"
$ws.Range('E34').Value = '  +6.94%  '
$ws.Range('D35').Formula = "'0.7660"
$ws.Range('E35').Value = '  +2.56%  '
$ws.Range('D36').Formula = "'2.978"
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').Formula = "'3.266"
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').Formula = "'2.593"
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Formula = "'0.5572"
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Formula = "'0.01981"
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').Formula = "'9.019"
$ws.Range('E42').Value = '  +5.87%  '
$ws.Range('D43').Formula = "'6.593"
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('D44').Formula = "'118.54"
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').Formula = "'0.1498"
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').Formula = "'0.4813"
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').Formula = "'0.9995"
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').Formula = "'10.14"
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('D50').Formula = "'37.62"
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('E51').Value = '  +1.18%  '
